# Append a new entry at the top of the data (row 2) on the "ランサーズ" sheet,
# pushing the existing rows down by one, and refresh the "取得日時" (fetched at)
# timestamp on every row to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-14 18:31:25"

# Snapshot existing rows 2..8 (A..H) before we start overwriting anything,
# so we can shift them down to rows 3..9.
# NOTE: use Value2 (not Value) to read cell contents back out reliably.
$firstRow = 2
$lastRow = 8

$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
    }
    $data += ,$row
}

# Clear out all existing hyperlinks up front; Range.Hyperlinks.Delete() in
# this environment operates on the whole sheet's hyperlink collection, so we
# only want to call it once, before re-adding hyperlinks in final order.
$ws.Range("A1").Hyperlinks.Delete()

# Write the shifted rows from the bottom up (old row N -> new row N+1) so we
# never clobber data before it has been read.
for ($i = $data.Length - 1; $i -ge 0; $i--) {
    $srcRow = $firstRow + $i
    $dstRow = $srcRow + 1
    $row = $data[$i]

    $ws.Cells.Item($dstRow, 1).Value2 = $newTimestamp
    $ws.Cells.Item($dstRow, 2).Value2 = $row.B
    $ws.Cells.Item($dstRow, 3).Value2 = $row.C
    $ws.Cells.Item($dstRow, 4).Value2 = $row.D
    $ws.Cells.Item($dstRow, 5).Value2 = $row.E

    $ws.Cells.Item($dstRow, 6).Value2 = $row.F
    $ws.Cells.Item($dstRow, 6).Style = "Hyperlink"

    $ws.Cells.Item($dstRow, 7).Value2 = $row.G

    if ($row.H) {
        $ws.Cells.Item($dstRow, 8).Value2 = $row.H
    } else {
        $ws.Cells.Item($dstRow, 8).ClearContents()
    }
}

# Now insert the brand-new entry into row 2.
$ws.Cells.Item(2, 1).Value2 = $newTimestamp
$ws.Cells.Item(2, 2).Value2 = "AIで精度の高い予想モデルやシステムを作って欲しいです。"
$ws.Cells.Item(2, 3).Value2 = "システム開発"
$ws.Cells.Item(2, 4).Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value2 = "期限情報なし"

$ws.Cells.Item(2, 6).Value2 = "https://www.lancers.jp/work/detail/5491912"
$ws.Cells.Item(2, 6).Style = "Hyperlink"

$ws.Cells.Item(2, 7).Value2 = 318
$ws.Cells.Item(2, 8).Value2 = "🔥AI,Ai"

# Re-create all hyperlinks (F2..F9) in order, now that every row is in its
# final place.
for ($r = 2; $r -le 9; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Range("F$r"), $target)
}
